$wb = $excel.ActiveWorkbook

# --- Neg_Change sheet ---
$ws1 = $wb.Worksheets.Item("Neg_Change")

# Row 2: RELIANCE
$ws1.Cells.Item(2, 1).Value = "RELIANCE"
$ws1.Cells.Item(2, 2).Value = 1465
$ws1.Cells.Item(2, 3).Value = 1479.9
$ws1.Cells.Item(2, 4).Value = 1465
$ws1.Cells.Item(2, 5).Value = 1476
$ws1.Cells.Item(2, 6).Value = 8335311
$ws1.Cells.Item(2, 7).Value = 16518684
$ws1.Cells.Item(2, 8).Value = -0.4954010258928617
$ws1.Cells.Item(2, 9).Value = "RELIANCE"

# Row 3: TCS
$ws1.Cells.Item(3, 1).Value = "TCS"
$ws1.Cells.Item(3, 2).Value = 3191
$ws1.Cells.Item(3, 3).Value = 3229.7
$ws1.Cells.Item(3, 4).Value = 3191
$ws1.Cells.Item(3, 5).Value = 3212.6
$ws1.Cells.Item(3, 6).Value = 1938673
$ws1.Cells.Item(3, 7).Value = 4159336
$ws1.Cells.Item(3, 8).Value = -0.5338984395586218
$ws1.Cells.Item(3, 9).Value = "TCS"

# Row 4: SBILIFE
$ws1.Cells.Item(4, 1).Value = "SBILIFE"
$ws1.Cells.Item(4, 2).Value = 2093.8
$ws1.Cells.Item(4, 3).Value = 2093.9
$ws1.Cells.Item(4, 4).Value = 2065
$ws1.Cells.Item(4, 5).Value = 2075
$ws1.Cells.Item(4, 6).Value = 455165
$ws1.Cells.Item(4, 7).Value = 956736
$ws1.Cells.Item(4, 8).Value = -0.5242522493143355
$ws1.Cells.Item(4, 9).Value = "SBILIFE"

# Row 5: APOLLOHOSP
$ws1.Cells.Item(5, 1).Value = "APOLLOHOSP"
$ws1.Cells.Item(5, 2).Value = 7348
$ws1.Cells.Item(5, 3).Value = 7360.5
$ws1.Cells.Item(5, 4).Value = 7225
$ws1.Cells.Item(5, 5).Value = 7255
$ws1.Cells.Item(5, 6).Value = 312461
$ws1.Cells.Item(5, 7).Value = 778407
$ws1.Cells.Item(5, 8).Value = -0.5985891699329529
$ws1.Cells.Item(5, 9).Value = "APOLLOHOSP"

# Row 6: BAJAJ-AUTO
$ws1.Cells.Item(6, 1).Value = "BAJAJ-AUTO"
$ws1.Cells.Item(6, 2).Value = 9760.5
$ws1.Cells.Item(6, 3).Value = 9830
$ws1.Cells.Item(6, 4).Value = 9538
$ws1.Cells.Item(6, 5).Value = 9578.5
$ws1.Cells.Item(6, 6).Value = 250994
$ws1.Cells.Item(6, 7).Value = 534086
$ws1.Cells.Item(6, 8).Value = -0.5300494676887243
$ws1.Cells.Item(6, 9).Value = "BAJAJ-AUTO"

# Row 7: ABB
$ws1.Cells.Item(7, 1).Value = "ABB"
$ws1.Cells.Item(7, 2).Value = 5050
$ws1.Cells.Item(7, 3).Value = 5138
$ws1.Cells.Item(7, 4).Value = 4983.5
$ws1.Cells.Item(7, 5).Value = 5074
$ws1.Cells.Item(7, 6).Value = 329338
$ws1.Cells.Item(7, 7).Value = 681840
$ws1.Cells.Item(7, 8).Value = -0.5169863897688607
$ws1.Cells.Item(7, 9).Value = "ABB"

# Row 8: ICICIPRULI
$ws1.Cells.Item(8, 1).Value = "ICICIPRULI"
$ws1.Cells.Item(8, 2).Value = 684
$ws1.Cells.Item(8, 3).Value = 693
$ws1.Cells.Item(8, 4).Value = 671.1
$ws1.Cells.Item(8, 5).Value = 687
$ws1.Cells.Item(8, 6).Value = 789798
$ws1.Cells.Item(8, 7).Value = 1909675
$ws1.Cells.Item(8, 8).Value = -0.5864228206370193
$ws1.Cells.Item(8, 9).Value = "ICICIPRULI"

# Row 9: DABUR
$ws1.Cells.Item(9, 1).Value = "DABUR"
$ws1.Cells.Item(9, 2).Value = 518.7
$ws1.Cells.Item(9, 3).Value = 525.55
$ws1.Cells.Item(9, 4).Value = 516.1
$ws1.Cells.Item(9, 5).Value = 523.7
$ws1.Cells.Item(9, 6).Value = 1522552
$ws1.Cells.Item(9, 7).Value = 3203359
$ws1.Cells.Item(9, 8).Value = -0.5247014149834596
$ws1.Cells.Item(9, 9).Value = "DABUR"

# Row 10: IDFCFIRSTB
$ws1.Cells.Item(10, 1).Value = "IDFCFIRSTB"
$ws1.Cells.Item(10, 2).Value = 85.90000000000001
$ws1.Cells.Item(10, 3).Value = 86.79000000000001
$ws1.Cells.Item(10, 4).Value = 85.06999999999999
$ws1.Cells.Item(10, 5).Value = 85.98
$ws1.Cells.Item(10, 6).Value = 31030611
$ws1.Cells.Item(10, 7).Value = 62702609
$ws1.Cells.Item(10, 8).Value = -0.5051145160482876
$ws1.Cells.Item(10, 9).Value = "IDFCFIRSTB"

# Row 11: POLICYBZR
$ws1.Cells.Item(11, 1).Value = "POLICYBZR"
$ws1.Cells.Item(11, 2).Value = 1689
$ws1.Cells.Item(11, 3).Value = 1700
$ws1.Cells.Item(11, 4).Value = 1668.1
$ws1.Cells.Item(11, 5).Value = 1687.8
$ws1.Cells.Item(11, 6).Value = 795693
$ws1.Cells.Item(11, 7).Value = 1665928
$ws1.Cells.Item(11, 8).Value = -0.5223725154988691
$ws1.Cells.Item(11, 9).Value = "POLICYBZR"

# Row 12: NYKAA
$ws1.Cells.Item(12, 1).Value = "NYKAA"
$ws1.Cells.Item(12, 2).Value = 257
$ws1.Cells.Item(12, 3).Value = 258.6
$ws1.Cells.Item(12, 4).Value = 253
$ws1.Cells.Item(12, 5).Value = 253.95
$ws1.Cells.Item(12, 6).Value = 3097714
$ws1.Cells.Item(12, 7).Value = 6101910
$ws1.Cells.Item(12, 8).Value = -0.4923369895655623
$ws1.Cells.Item(12, 9).Value = "NYKAA"

# Row 13: KALYANKJIL
$ws1.Cells.Item(13, 1).Value = "KALYANKJIL"
$ws1.Cells.Item(13, 2).Value = 509.3
$ws1.Cells.Item(13, 3).Value = 515.45
$ws1.Cells.Item(13, 4).Value = 502.75
$ws1.Cells.Item(13, 5).Value = 505.95
$ws1.Cells.Item(13, 6).Value = 2355831
$ws1.Cells.Item(13, 7).Value = 5202626
$ws1.Cells.Item(13, 8).Value = -0.5471842488773938
$ws1.Cells.Item(13, 9).Value = "KALYANKJIL"

# Row 14: TATATECH
$ws1.Cells.Item(14, 1).Value = "TATATECH"
$ws1.Cells.Item(14, 2).Value = 660.5
$ws1.Cells.Item(14, 3).Value = 664.85
$ws1.Cells.Item(14, 4).Value = 648.7
$ws1.Cells.Item(14, 5).Value = 652.25
$ws1.Cells.Item(14, 6).Value = 903801
$ws1.Cells.Item(14, 7).Value = 1864357
$ws1.Cells.Item(14, 8).Value = -0.5152210654933578
$ws1.Cells.Item(14, 9).Value = "TATATECH"

# Row 15: IRB
$ws1.Cells.Item(15, 1).Value = "IRB"
$ws1.Cells.Item(15, 2).Value = 41.76
$ws1.Cells.Item(15, 3).Value = 42.08
$ws1.Cells.Item(15, 4).Value = 41.32
$ws1.Cells.Item(15, 5).Value = 41.48
$ws1.Cells.Item(15, 6).Value = 7252215
$ws1.Cells.Item(15, 7).Value = 17749470
$ws1.Cells.Item(15, 8).Value = -0.5914123069590247
$ws1.Cells.Item(15, 9).Value = "IRB"

# Row 16: TORNTPOWER
$ws1.Cells.Item(16, 1).Value = "TORNTPOWER"
$ws1.Cells.Item(16, 2).Value = 1350
$ws1.Cells.Item(16, 3).Value = 1367.2
$ws1.Cells.Item(16, 4).Value = 1318.5
$ws1.Cells.Item(16, 5).Value = 1320.1
$ws1.Cells.Item(16, 6).Value = 331050
$ws1.Cells.Item(16, 7).Value = 774422
$ws1.Cells.Item(16, 8).Value = -0.572519892255127
$ws1.Cells.Item(16, 9).Value = "TORNTPOWER"

# --- Pos_Change sheet ---
$ws2 = $wb.Worksheets.Item("Pos_Change")

# Row 2: ETERNAL
$ws2.Cells.Item(2, 1).Value = "ETERNAL"
$ws2.Cells.Item(2, 2).Value = 283.4
$ws2.Cells.Item(2, 3).Value = 292.9
$ws2.Cells.Item(2, 4).Value = 281.65
$ws2.Cells.Item(2, 5).Value = 284.7
$ws2.Cells.Item(2, 6).Value = 55851766
$ws2.Cells.Item(2, 7).Value = 38191591
$ws2.Cells.Item(2, 8).Value = 0.4624100368062697
$ws2.Cells.Item(2, 9).Value = "ETERNAL"

# Row 3: AXISBANK
$ws2.Cells.Item(3, 1).Value = "AXISBANK"
$ws2.Cells.Item(3, 2).Value = 1286.8
$ws2.Cells.Item(3, 3).Value = 1295.3
$ws2.Cells.Item(3, 4).Value = 1270.1
$ws2.Cells.Item(3, 5).Value = 1274.5
$ws2.Cells.Item(3, 6).Value = 5531427
$ws2.Cells.Item(3, 7).Value = 3949654
$ws2.Cells.Item(3, 8).Value = 0.4004839411249694
$ws2.Cells.Item(3, 9).Value = "AXISBANK"

# Row 4: BAJAJFINSV
$ws2.Cells.Item(4, 1).Value = "BAJAJFINSV"
$ws2.Cells.Item(4, 2).Value = 2008.9
$ws2.Cells.Item(4, 3).Value = 2014.2
$ws2.Cells.Item(4, 4).Value = 1981.2
$ws2.Cells.Item(4, 5).Value = 1989.1
$ws2.Cells.Item(4, 6).Value = 1143760
$ws2.Cells.Item(4, 7).Value = 756940
$ws2.Cells.Item(4, 8).Value = 0.5110312574312363
$ws2.Cells.Item(4, 9).Value = "BAJAJFINSV"

# Row 5: DMART
$ws2.Cells.Item(5, 1).Value = "DMART"
$ws2.Cells.Item(5, 2).Value = 3794
$ws2.Cells.Item(5, 3).Value = 3833.2
$ws2.Cells.Item(5, 4).Value = 3745.1
$ws2.Cells.Item(5, 5).Value = 3807
$ws2.Cells.Item(5, 6).Value = 600885
$ws2.Cells.Item(5, 7).Value = 400228
$ws2.Cells.Item(5, 8).Value = 0.5013567266658004
$ws2.Cells.Item(5, 9).Value = "DMART"

# Row 6: TVSMOTOR
$ws2.Cells.Item(6, 1).Value = "TVSMOTOR"
$ws2.Cells.Item(6, 2).Value = 3800
$ws2.Cells.Item(6, 3).Value = 3848.7
$ws2.Cells.Item(6, 4).Value = 3751.3
$ws2.Cells.Item(6, 5).Value = 3775.9
$ws2.Cells.Item(6, 6).Value = 567897
$ws2.Cells.Item(6, 7).Value = 367668
$ws2.Cells.Item(6, 8).Value = 0.5445918600476517
$ws2.Cells.Item(6, 9).Value = "TVSMOTOR"

# Row 7: TORNTPHARM
$ws2.Cells.Item(7, 1).Value = "TORNTPHARM"
$ws2.Cells.Item(7, 2).Value = 4016
$ws2.Cells.Item(7, 3).Value = 4037.6
$ws2.Cells.Item(7, 4).Value = 3950
$ws2.Cells.Item(7, 5).Value = 3955
$ws2.Cells.Item(7, 6).Value = 312744
$ws2.Cells.Item(7, 7).Value = 202187
$ws2.Cells.Item(7, 8).Value = 0.5468056798903985
$ws2.Cells.Item(7, 9).Value = "TORNTPHARM"

# Row 8: LODHA
$ws2.Cells.Item(8, 1).Value = "LODHA"
$ws2.Cells.Item(8, 2).Value = 1082
$ws2.Cells.Item(8, 3).Value = 1082.1
$ws2.Cells.Item(8, 4).Value = 1052.8
$ws2.Cells.Item(8, 5).Value = 1064
$ws2.Cells.Item(8, 6).Value = 2660754
$ws2.Cells.Item(8, 7).Value = 1683229
$ws2.Cells.Item(8, 8).Value = 0.5807439154149554
$ws2.Cells.Item(8, 9).Value = "LODHA"

# Row 9: TATAPOWER
$ws2.Cells.Item(9, 1).Value = "TATAPOWER"
$ws2.Cells.Item(9, 2).Value = 373
$ws2.Cells.Item(9, 3).Value = 374
$ws2.Cells.Item(9, 4).Value = 362.95
$ws2.Cells.Item(9, 5).Value = 364.3
$ws2.Cells.Item(9, 6).Value = 5875179
$ws2.Cells.Item(9, 7).Value = 3791372
$ws2.Cells.Item(9, 8).Value = 0.5496181857121907
$ws2.Cells.Item(9, 9).Value = "TATAPOWER"

# Row 10: IGL
$ws2.Cells.Item(10, 1).Value = "IGL"
$ws2.Cells.Item(10, 2).Value = 185.6
$ws2.Cells.Item(10, 3).Value = 187.95
$ws2.Cells.Item(10, 4).Value = 183
$ws2.Cells.Item(10, 5).Value = 187.35
$ws2.Cells.Item(10, 6).Value = 1249424
$ws2.Cells.Item(10, 7).Value = 831626
$ws2.Cells.Item(10, 8).Value = 0.5023868902607662
$ws2.Cells.Item(10, 9).Value = "IGL"

# Row 11: OBEROIRLTY
$ws2.Cells.Item(11, 1).Value = "OBEROIRLTY"
$ws2.Cells.Item(11, 2).Value = 1725.5
$ws2.Cells.Item(11, 3).Value = 1732.2
$ws2.Cells.Item(11, 4).Value = 1665.3
$ws2.Cells.Item(11, 5).Value = 1699
$ws2.Cells.Item(11, 6).Value = 857630
$ws2.Cells.Item(11, 7).Value = 564739
$ws2.Cells.Item(11, 8).Value = 0.5186307303019625
$ws2.Cells.Item(11, 9).Value = "OBEROIRLTY"

# Row 12: ABCAPITAL
$ws2.Cells.Item(12, 1).Value = "ABCAPITAL"
$ws2.Cells.Item(12, 2).Value = 357
$ws2.Cells.Item(12, 3).Value = 363.3
$ws2.Cells.Item(12, 4).Value = 351.45
$ws2.Cells.Item(12, 5).Value = 353.7
$ws2.Cells.Item(12, 6).Value = 3959855
$ws2.Cells.Item(12, 7).Value = 2771525
$ws2.Cells.Item(12, 8).Value = 0.4287639476461515
$ws2.Cells.Item(12, 9).Value = "ABCAPITAL"

# Row 13: MANKIND
$ws2.Cells.Item(13, 1).Value = "MANKIND"
$ws2.Cells.Item(13, 2).Value = 2250
$ws2.Cells.Item(13, 3).Value = 2250.7
$ws2.Cells.Item(13, 4).Value = 2186
$ws2.Cells.Item(13, 5).Value = 2201.3
$ws2.Cells.Item(13, 6).Value = 396485
$ws2.Cells.Item(13, 7).Value = 271109
$ws2.Cells.Item(13, 8).Value = 0.4624560601086648
$ws2.Cells.Item(13, 9).Value = "MANKIND"

# Row 14: IREDA
$ws2.Cells.Item(14, 1).Value = "IREDA"
$ws2.Cells.Item(14, 2).Value = 141
$ws2.Cells.Item(14, 3).Value = 141.64
$ws2.Cells.Item(14, 4).Value = 135.6
$ws2.Cells.Item(14, 5).Value = 137.1
$ws2.Cells.Item(14, 6).Value = 15649670
$ws2.Cells.Item(14, 7).Value = 10554455
$ws2.Cells.Item(14, 8).Value = 0.4827549124990348
$ws2.Cells.Item(14, 9).Value = "IREDA"

# Row 15: SUZLON
$ws2.Cells.Item(15, 1).Value = "SUZLON"
$ws2.Cells.Item(15, 2).Value = 50.82
$ws2.Cells.Item(15, 3).Value = 51.09
$ws2.Cells.Item(15, 4).Value = 48.87
$ws2.Cells.Item(15, 5).Value = 49.3
$ws2.Cells.Item(15, 6).Value = 89917429
$ws2.Cells.Item(15, 7).Value = 59316617
$ws2.Cells.Item(15, 8).Value = 0.5158893670554408
$ws2.Cells.Item(15, 9).Value = "SUZLON"

# Row 16: IIFL
$ws2.Cells.Item(16, 1).Value = "IIFL"
$ws2.Cells.Item(16, 2).Value = 647.9
$ws2.Cells.Item(16, 3).Value = 656.7
$ws2.Cells.Item(16, 4).Value = 638.65
$ws2.Cells.Item(16, 5).Value = 652.3
$ws2.Cells.Item(16, 6).Value = 1846485
$ws2.Cells.Item(16, 7).Value = 1240235
$ws2.Cells.Item(16, 8).Value = 0.4888186513039867
$ws2.Cells.Item(16, 9).Value = "IIFL"

# Row 17: CROMPTON
$ws2.Cells.Item(17, 1).Value = "CROMPTON"
$ws2.Cells.Item(17, 2).Value = 257.75
$ws2.Cells.Item(17, 3).Value = 259.05
$ws2.Cells.Item(17, 4).Value = 250
$ws2.Cells.Item(17, 5).Value = 251.75
$ws2.Cells.Item(17, 6).Value = 3912266
$ws2.Cells.Item(17, 7).Value = 2677012
$ws2.Cells.Item(17, 8).Value = 0.4614301318036677
$ws2.Cells.Item(17, 9).Value = "CROMPTON"

# Row 18: CAMS
$ws2.Cells.Item(18, 1).Value = "CAMS"
$ws2.Cells.Item(18, 2).Value = 744.7
$ws2.Cells.Item(18, 3).Value = 747.55
$ws2.Cells.Item(18, 4).Value = 721.05
$ws2.Cells.Item(18, 5).Value = 727.4
$ws2.Cells.Item(18, 6).Value = 1229863
$ws2.Cells.Item(18, 7).Value = 844354
$ws2.Cells.Item(18, 8).Value = 0.4565727171304926
$ws2.Cells.Item(18, 9).Value = "CAMS"

